# Clean up the "roboticS1Prep" column (I): replace the text "No" with the
# boolean FALSE value, displayed with a custom TRUE/FALSE number format,
# for every data row (rows 2-37).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("I2:I37")
$rng.Value = $false
$rng.NumberFormat = '"TRUE";"TRUE";"FALSE"'

# Match the author's final selection/active cell (now on column I instead
# of the old column H selection).
$rng.Select()
